# Update "想去人数" (interested-count) figures that changed in the upstream
# data refresh, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 330
    $ws.Range("F5").Value = 4884
}
